$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B2" = 0.1909620991253644
    "C2" = 0.5626822157434402
    "J2" = 0.01457725947521866
    "P2" = 0.1516034985422741
    "S2" = 0.08017492711370262
    "B3" = 0.004830917874396135
    "C3" = 0.01932367149758454
    "J3" = 0.02173913043478261
    "P3" = 0.7971014492753623
    "S3" = 0.1570048309178744
    "J4" = 0.05833333333333333
    "P4" = 0.725
    "S4" = 0.2166666666666667
    "B6" = 0.06326034063260341
    "D6" = 0.0194647201946472
    "E6" = 0.0024330900243309
    "F6" = 0.072992700729927
    "J6" = 0.2287104622871046
    "O6" = 0.024330900243309
    "Q6" = 0.218978102189781
    "R6" = 0.0583941605839416
    "S6" = 0.3114355231143552
    "B7" = 0.1233766233766234
    "D7" = 0.0303030303030303
    "E7" = 0.006493506493506494
    "F7" = 0.0670995670995671
    "J7" = 0.1428571428571428
    "O7" = 0.02164502164502164
    "Q7" = 0.1731601731601732
    "R7" = 0.08441558441558442
    "S7" = 0.3506493506493507
    "B8" = 0.1052631578947368
    "D8" = 0.02799552071668533
    "E8" = 0.001119820828667413
    "F8" = 0.05375139977603583
    "J8" = 0.1175811870100784
    "O8" = 0.02015677491601344
    "Q8" = 0.2138857782754759
    "R8" = 0.05487122060470324
    "S8" = 0.4053751399776036
    "B9" = 0.1157407407407407
    "D9" = 0.02546296296296296
    "E9" = 0.002314814814814815
    "F9" = 0.03703703703703703
    "J9" = 0.1388888888888889
    "O9" = 0.0162037037037037
    "Q9" = 0.1944444444444444
    "R9" = 0.07175925925925926
    "S9" = 0.3981481481481481
    "B10" = 0.1214340786430224
    "D10" = 0.02467232074016962
    "E10" = 0.001156515034695451
    "F10" = 0.06630686198920586
    "J10" = 0.1283731688511951
    "O10" = 0.01117964533538936
    "Q10" = 0.2197378565921357
    "R10" = 0.07787201233616037
    "S10" = 0.3492675404780262
    "F11" = 0.001358695652173913
    "G11" = 0.1290760869565217
    "J11" = 0.08831521739130435
    "K11" = 0.1888586956521739
    "L11" = 0.5665760869565217
    "S11" = 0.02581521739130435
    "G12" = 0.7428571428571429
    "J12" = 0.2
    "K12" = 0.007142857142857143
    "L12" = 0.02142857142857143
    "S12" = 0.02857142857142857
    "G13" = 0.6666666666666666
    "J13" = 0.2843137254901961
    "S13" = 0.04901960784313725
    "G14" = 0.8333333333333334
    "J14" = 0.1666666666666667
    "F15" = 0.01869158878504673
    "H15" = 0.1425233644859813
    "I15" = 0.05373831775700934
    "J15" = 0.352803738317757
    "K15" = 0.1051401869158878
    "M15" = 0.02102803738317757
    "O15" = 0.05841121495327103
    "S15" = 0.2476635514018692
    "F16" = 0.01151631477927063
    "H16" = 0.1900191938579655
    "I16" = 0.08445297504798464
    "J16" = 0.3704414587332054
    "K16" = 0.1381957773512476
    "M16" = 0.02687140115163148
    "N16" = 0.003838771593090211
    "O16" = 0.04606525911708254
    "S16" = 0.1285988483685221
    "F17" = 0.01081612586037365
    "H17" = 0.1553588987217306
    "I17" = 0.09734513274336283
    "J17" = 0.4336283185840708
    "K17" = 0.119960668633235
    "M17" = 0.01966568338249754
    "O17" = 0.04621435594886922
    "S17" = 0.1170108161258604
    "F18" = 0.01146131805157593
    "H18" = 0.1805157593123209
    "I18" = 0.09742120343839542
    "J18" = 0.3782234957020058
    "K18" = 0.1088825214899714
    "M18" = 0.02865329512893983
    "O18" = 0.07736389684813753
    "S18" = 0.1174785100286533
    "F19" = 0.01569416498993964
    "H19" = 0.2056338028169014
    "I19" = 0.09255533199195171
    "J19" = 0.3541247484909457
    "K19" = 0.1251509054325956
    "M19" = 0.02293762575452716
    "N19" = 0.001609657947686117
    "O19" = 0.06841046277665996
    "S19" = 0.1138832997987928
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
